# Commit: "added poisson calculation in main.db"
# Fills in the previously-blank probability inputs (columns A/B on sheet
# "2way" and A/B/C on sheet "3way") for rows 187-195, which drive the
# dependent odds/value formulas that currently evaluate to #DIV/0!.
# Also fills in the predicted-score columns (T/U) on sheet "2way" for the
# same rows.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "2way": probabilities (A,B) and predicted score (T,U)
# ------------------------------------------------------------------
$ws2way = $wb.Worksheets.Item("2way")

$twoWayData = @{
    187 = @(0.22, 0.78, "X", "0-0")
    188 = @(0.8,  0.2,  "2", "1-2")
    189 = @(0.48, 0.52, "1", "1-0")
    190 = @(0.56000000000000005, 0.44, "1", "2-0")
    191 = @(0.53, 0.47, "1", "2-0")
    192 = @(0.36, 0.64, "1", "1-0")
    193 = @(0.64, 0.36, "1", "2-1")
    194 = @(0.39, 0.61, "1", "1-0")
    195 = @(0.44, 0.56000000000000005, "2", "0-1")
}

foreach ($row in 187..195) {
    $vals = $twoWayData[$row]
    $ws2way.Cells.Item($row, 1).Value = $vals[0]   # A = Home probability
    $ws2way.Cells.Item($row, 2).Value = $vals[1]   # B = Away probability
    $ws2way.Cells.Item($row, 20).Value = $vals[2]  # T = predicted result
    $ws2way.Cells.Item($row, 21).Value = $vals[3]  # U = predicted score
}

# ------------------------------------------------------------------
# Sheet "3way": probabilities (A,B,C)
# ------------------------------------------------------------------
$ws3way = $wb.Worksheets.Item("3way")

$threeWayData = @{
    187 = @(0.37129934086969912, 0.34884400310672858, 0.2798566525160639)
    188 = @(0.23463255043551384, 0.18073184744129223, 0.58453651640354176)
    189 = @(0.64022647314302727, 0.21705378938130396, 0.14271536689544242)
    190 = @(0.69437545673915368, 0.18693221088767259, 0.11867617256818047)
    191 = @(0.78569468180377411, 0.15203477560086567, 0.06224049241555441)
    192 = @(0.60094847821193087, 0.2552177012388635,  0.14383314036308029)
    193 = @(0.47791523624913868, 0.22533944083901677, 0.29674073680652535)
    194 = @(0.70368336387766373, 0.21055529732355346, 0.08575831802525442)
    195 = @(0.14233193326741073, 0.22849422328188518, 0.62917141647687258)
}

foreach ($row in 187..195) {
    $vals = $threeWayData[$row]
    $ws3way.Cells.Item($row, 1).Value = $vals[0]   # A = Home probability
    $ws3way.Cells.Item($row, 2).Value = $vals[1]   # B = Draw probability
    $ws3way.Cells.Item($row, 3).Value = $vals[2]   # C = Away probability
}

# ------------------------------------------------------------------
# Restore the view/selection state left behind after the edit
# ------------------------------------------------------------------
$ws2way.Range("U196").Select()
$ws3way.Range("A195:C195").Select()
$ws2way.Activate()
